# Onboarding-bootstrap datav4: add an "id" column to the "locations" sheet
# and make "locations" the active sheet.

$wb = $excel.ActiveWorkbook
$wsLoc = $wb.Worksheets.Item("locations")

# Insert a new column A before the existing data, shifting hotels' columns
# (and their formatting/styles) one place to the right.
$wsLoc.Columns.Item(1).Insert() | Out-Null

# Populate the new id column: header + numeric ids 200000-200008.
$wsLoc.Range("A1").Value = "id"
$ids = @(200000, 200001, 200002, 200003, 200004, 200005, 200006, 200007, 200008)
for ($i = 0; $i -lt $ids.Length; $i++) {
    $wsLoc.Cells.Item($i + 2, 1).Value = $ids[$i]
}

# "hotels" keeps a plain view (no tabSelected/topLeftCell anymore) but the
# same selection shape as before (A2 anchored, full data column selected).
$wsHotels = $wb.Worksheets.Item("hotels")
$wsHotels.Activate() | Out-Null
$wsHotels.Range("A2:A107").Select() | Out-Null

# Make "locations" the active sheet/tab, with the same selection shape the
# "hotels" sheet used to have (A2 anchored, full data column selected).
$wsLoc.Activate() | Out-Null
$wsLoc.Range("A2:A10").Select() | Out-Null

# Update the "locations" defined name to reflect the shifted columns.
$wb.Names.Item("locations").RefersTo = "=locations!`$B`$1:`$K`$10"

Write-Host "done"
